# Daily attendance processing - 2025-12-08 01:28:01
# Normalizes the "Recorded By" column (G) so that the most recently
# recording user (the last entry in the comma-separated list) is moved
# to the front of the list - except when that last entry is the
# generic "System"/"system" placeholder, which is left in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if ($value -notlike "*,*") { continue }

    $parts = $value -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $n = $parts.Length
    $lastPart = $parts[$n - 1]
    if ($lastPart.ToLower() -eq "system") {
        continue
    }

    $newParts = @($lastPart) + $parts[0..($n - 2)]
    $newValue = [string]::Join(", ", $newParts)

    $cell.Value2 = $newValue
}
